$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 42.5
$ws.Columns.Item(3).ColumnWidth = 69.83333333333333
$ws.Columns.Item(4).ColumnWidth = 40.166666666666664

# --- Cell values + per-cell formatting (8pt font; wrap where needed) ---
$ws.Range("A2").Font.Size = 8
$ws.Range("A2").Value = "JDK"

$ws.Range("B2").Font.Size = 8
$ws.Range("B2").WrapText = $true
$ws.Range("B2").Value = "java.util.Properties`njava.util.Hashtable`njava.util.Enumeration<T>"

$ws.Range("C2").Font.Size = 8
$ws.Range("C2").WrapText = $true
$ws.Range("C2").Value = "一個簡單的接口, hasMoreElement和nextElement, 在Properties看到有所應用, 當要遍歷HashTable時, 將h,key()賦給一個臨時的Enumeration 類(這key大有文章, 是由HashTable中的inner class Enumerator 實現的), 再以for 的三段式利用hasMoreElement()和nextElement()來Iterate"

$ws.Range("D2").Font.Size = 8
$ws.Range("D2").WrapText = $true
$ws.Range("D2").Value = "It is fun, like music, like literature, like art"

$ws.Range("A3").Font.Size = 8
$ws.Range("A3").Value = "Concurrency"

$ws.Range("B3").Font.Size = 8
$ws.Range("B3").Value = "java.lang.Thread"

$ws.Range("C3").Font.Size = 8
$ws.Range("C3").WrapText = $true
$ws.Range("C3").Value = "synchronized method in Thread"

$ws.Range("D3").Font.Size = 8
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Value = "the beauty in mathematic"

$ws.Range("A4").Font.Size = 8
$ws.Range("A4").Value = "Collection"

$ws.Range("B4").Font.Size = 8
$ws.Range("B4").Value = "java.util.Collection"

$ws.Range("C4").Font.Size = 8
$ws.Range("C4").WrapText = $true
$ws.Range("C4").Value = "keyword `"default`" is added in jdk8 to extend interfaces to adopt lambda expression. So that adding new method in interface wont break the compilation and running of java code in new jre. E.g java.util.Collection"

$ws.Range("D4").Font.Size = 8
$ws.Range("D4").WrapText = $true
$ws.Range("D4").Value = "the art of code, the message between lines"

$ws.Range("A5").Font.Size = 8
$ws.Range("A5").Value = "JDK"

$ws.Range("B5").Font.Size = 8
$ws.Range("B5").Value = "java.lang.Boolean"

$ws.Range("C5").Font.Size = 8
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Value = "read java.lang.Boolean"

$ws.Range("A6").Font.Size = 8
$ws.Range("A6").Value = "JDK"

$ws.Range("B6").Font.Size = 8
$ws.Range("B6").WrapText = $true
$ws.Range("B6").Value = "java.lang.System"

$ws.Range("C6").Font.Size = 8
$ws.Range("C6").WrapText = $true
$ws.Range("C6").Value = "File intputing for java utilities`nUse System.getProperty(`"user.dir`") can empower java utilities to take file input from relative filepath"

$ws.Range("A7").Font.Size = 8
$ws.Range("A7").Value = "General"

$ws.Range("B7").Font.Size = 8
$ws.Range("B7").Value = "Use Eclipse to read code"

$ws.Range("C7").Font.Size = 8
$ws.Range("C7").Value = "<detail refers to index.xlsx where I note how to setup maven source sync env in eclipse>"

$ws.Range("A8").Font.Size = 8
$ws.Range("A8").Value = "Jtopen"

$ws.Range("B8").Font.Size = 8
$ws.Range("B8").WrapText = $true
$ws.Range("B8").Value = "com.ibm.as400.access.AS400`ncom.ibm.as400.access.ProgramParameter`ncom.ibm.as400.access.ProgramCall"

$ws.Range("C8").Font.Size = 8
$ws.Range("C8").WrapText = $true
$ws.Range("C8").Value = "* Learn that jtopen is more flat in terms of structure`n* it used sync/ transient a lot`n* its private member is named like name_`n* programparameter's core is byte[] for input / output data, and int for length`n* In AS400 class, it use signon(args…) to form a routing of the same method`n    signon()  > signon(Sting, String) > signon(String, String, Credit)...`n* The signon behind AS400 is performed by AS400Impl interface, and this interface is implemented by AS400ImplProxy and AS400ImplRemote. `n* Above interface and implementation relations is common in jtopen like the core behind ProgramCall : RemoteCommandImpl and its implementers Proxy and Remote"

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 45.75
$ws.Rows.Item(4).RowHeight = 34.5
$ws.Rows.Item(6).RowHeight = 23.25
$ws.Rows.Item(8).RowHeight = 113.25

# --- Selection ---
$ws.Range("C9").Select() | Out-Null
